# The sheet contained a final "Total Geral" row (30) followed by a
# "Desconsiderando ..." summary row (31). The "Total Geral" row is removed,
# so the "Desconsiderando ..." row shifts up to become row 30 and the used
# range shrinks from A1:N31 to A1:N30.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(30).Delete()
